$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data regenerated to filter save games.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win, G=sum
$data = @{
    2  = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    3  = @(3.286832544864788,  0.306821227259698,   0.1494219747398047, 0.4942365360607697, 4.23731228292506)
    4  = @(1.455362044514542,  1.655778082260271,   3.537761648806719,  0.4942365360607697, 7.143138311642302)
    5  = @(0.6606524410359556, 1.655778082260271,   22.3905356188092,   10.19245300693656,  34.89941914904198)
    6  = @(0.1190320826869504, 0.306821227259698,   0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    7  = @(3.286832544864788,  1.655778082260271,   0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    8  = @(3.286832544864788,  1.655778082260271,   0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    9  = @(0.04271373187048222,1.655778082260271,   0.1494219747398047, 0.4942365360607697, 2.342150324931327)
    10 = @(3.286832544864788,  1.655778082260271,   0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    11 = @(3.286832544864788,  1.655778082260271,   0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    12 = @(3.286832544864788,  1.655778082260271,   0.1494219747398047, 0.4942365360607697, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
